# Add Jan 2018 data to the journal URL tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing F242 URL cell was missing its hyperlink; add it back,
# copying the hyperlink-style formatting already used elsewhere in column F.
$ws.Hyperlinks.Add($ws.Range("F242"), "https://academic.oup.com/ps/article/97/1/131/4565714")
$ws.Range("F242").Style = $ws.Range("F2").Style

# Append new January 2018 "journal of animal science" rows (257-272):
# monogastric articles first, then ruminant articles.
$ws.Range("A257").Value = "journal of animal science"
$ws.Range("B257").Value = 2018
$ws.Range("C257").Value = 1
$ws.Range("D257").Value = "monogastric"
$ws.Range("E257").Value = "original article"
$ws.Range("F257").Value = "https://academic.oup.com/jas/article/96/1/98/4844089"

$ws.Range("A258").Value = "journal of animal science"
$ws.Range("B258").Value = 2018
$ws.Range("C258").Value = 1
$ws.Range("D258").Value = "monogastric"
$ws.Range("E258").Value = "original article"
$ws.Range("F258").Value = "https://academic.oup.com/jas/article/96/1/168/4844098"

$ws.Range("A259").Value = "journal of animal science"
$ws.Range("B259").Value = 2018
$ws.Range("C259").Value = 1
$ws.Range("D259").Value = "monogastric"
$ws.Range("E259").Value = "original article"
$ws.Range("F259").Value = "https://academic.oup.com/jas/article/96/1/181/4827713"

$ws.Range("A260").Value = "journal of animal science"
$ws.Range("B260").Value = 2018
$ws.Range("C260").Value = 1
$ws.Range("D260").Value = "monogastric"
$ws.Range("E260").Value = "original article"
$ws.Range("F260").Value = "https://academic.oup.com/jas/article/96/1/194/4827625"

$ws.Range("A261").Value = "journal of animal science"
$ws.Range("B261").Value = 2018
$ws.Range("C261").Value = 1
$ws.Range("D261").Value = "monogastric"
$ws.Range("E261").Value = "original article"
$ws.Range("F261").Value = "https://academic.oup.com/jas/article/96/1/206/4824921"

$ws.Range("A262").Value = "journal of animal science"
$ws.Range("B262").Value = 2018
$ws.Range("C262").Value = 1
$ws.Range("D262").Value = "monogastric"
$ws.Range("E262").Value = "original article"
$ws.Range("F262").Value = "https://academic.oup.com/jas/article/96/1/215/4844080"

$ws.Range("A263").Value = "journal of animal science"
$ws.Range("B263").Value = 2018
$ws.Range("C263").Value = 1
$ws.Range("D263").Value = "monogastric"
$ws.Range("E263").Value = "original article"
$ws.Range("F263").Value = "https://academic.oup.com/jas/article/96/1/225/4824872"

$ws.Range("A264").Value = "journal of animal science"
$ws.Range("B264").Value = 2018
$ws.Range("C264").Value = 1
$ws.Range("D264").Value = "ruminant"
$ws.Range("E264").Value = "original article"
$ws.Range("F264").Value = "https://academic.oup.com/jas/article/96/1/273/4827718"

$ws.Range("A265").Value = "journal of animal science"
$ws.Range("B265").Value = 2018
$ws.Range("C265").Value = 1
$ws.Range("D265").Value = "ruminant"
$ws.Range("E265").Value = "original article"
$ws.Range("F265").Value = "https://academic.oup.com/jas/article/96/1/284/4827785"

$ws.Range("A266").Value = "journal of animal science"
$ws.Range("B266").Value = 2018
$ws.Range("C266").Value = 1
$ws.Range("D266").Value = "ruminant"
$ws.Range("E266").Value = "original article"
$ws.Range("F266").Value = "https://academic.oup.com/jas/article/96/1/293/4827629"

$ws.Range("A267").Value = "journal of animal science"
$ws.Range("B267").Value = 2018
$ws.Range("C267").Value = 1
$ws.Range("D267").Value = "ruminant"
$ws.Range("E267").Value = "original article"
$ws.Range("F267").Value = "https://academic.oup.com/jas/article/96/1/306/4844085"

$ws.Range("A268").Value = "journal of animal science"
$ws.Range("B268").Value = 2018
$ws.Range("C268").Value = 1
$ws.Range("D268").Value = "ruminant"
$ws.Range("E268").Value = "original article"
$ws.Range("F268").Value = "https://academic.oup.com/jas/article/96/1/318/4825262"

$ws.Range("A269").Value = "journal of animal science"
$ws.Range("B269").Value = 2018
$ws.Range("C269").Value = 1
$ws.Range("D269").Value = "ruminant"
$ws.Range("E269").Value = "original article"
$ws.Range("F269").Value = "https://academic.oup.com/jas/article/96/1/331/4825179"

$ws.Range("A270").Value = "journal of animal science"
$ws.Range("B270").Value = 2018
$ws.Range("C270").Value = 1
$ws.Range("D270").Value = "ruminant"
$ws.Range("E270").Value = "original article"
$ws.Range("F270").Value = "https://academic.oup.com/jas/article/96/1/343/4818648"

$ws.Range("A271").Value = "journal of animal science"
$ws.Range("B271").Value = 2018
$ws.Range("C271").Value = 1
$ws.Range("D271").Value = "ruminant"
$ws.Range("E271").Value = "original article"
$ws.Range("F271").Value = "https://academic.oup.com/jas/article/96/1/354/4818673"

$ws.Range("A272").Value = "journal of animal science"
$ws.Range("B272").Value = 2018
$ws.Range("C272").Value = 1
$ws.Range("D272").Value = "ruminant"
$ws.Range("E272").Value = "original article"
$ws.Range("F272").Value = "https://academic.oup.com/jas/article/96/1/364/4818649"

# Leave the selection where the author ended up after entering the new rows.
$ws.Range("D279").Select()
